$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.563.57"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.958.24"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'244.32"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").Value = "'58.67"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.378"
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'0.0811"
$ws.Range("E10").Value = "  -5.43%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'22.13"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'0.830"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "2.245.63"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'13.72"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "1.953.18"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "36.498.24"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'69.73"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "'228.56"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "'5.06"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'9.24"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "'0.140"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'160.36"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").Value = "'19.43"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "'0.0621"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("D34").Value = "'4.32"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'2.24"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "'3.39"
$ws.Range("E37").Value = "  +10.50%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'5.75"
$ws.Range("E39").Value = "  -10.75%  "
$ws.Range("D40").Value = "'0.0981"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'1.18"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "'0.0212"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "'16.04"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "1.366.78"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'87.92"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "'7.15"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "2.136.57"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'43.71"
$ws.Range("E51").Value = "  -5.24%  "
